$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ActorTable")

# New headers
$ws.Cells.Item(1, 8).Value = "attackRange|Float"
$ws.Cells.Item(1, 9).Value = "ultimateRange|Float"

# attackRange values for rows 2..18 (row 3 = 4, row 4 = 2, rest = 0)
$attackRange = @{
    2  = 0
    3  = 4
    4  = 2
    5  = 0
    6  = 0
    7  = 0
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
}

for ($r = 2; $r -le 18; $r++) {
    $ws.Cells.Item($r, 8).Value = $attackRange[$r]
    $ws.Cells.Item($r, 9).Value = 0
}
